# Generate Report for Handback
# Updates timestamps (and one status code) in the handback-status workbook
# to reflect a later report-generation run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" column (G) on the Overview sheet, rows 2-3
$overview.Range("G2").Value = "2016-08-16 18:13:33"
$overview.Range("G3").Value = "2016-08-16 18:13:33"

# Priority column (E) on zh-cn, rows 2-3: "ht" -> "mt"
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E3").Value = "mt"

# Correspond Handoff Datetime (H) on zh-cn, rows 2-3
$zhcn.Range("H2").Value = "2016-08-16 18:13:28"
$zhcn.Range("H3").Value = "2016-08-16 18:13:28"

# Correspond Handback DateTime (K) on zh-cn, rows 2-3
$zhcn.Range("K2").Value = "2016-08-16 18:13:46"
$zhcn.Range("K3").Value = "2016-08-16 18:13:46"

# Priority column (E) on de-de, rows 2-3: "ht" -> "mt"
$dede.Range("E2").Value = "mt"
$dede.Range("E3").Value = "mt"

# Correspond Handoff Datetime (H) on de-de, rows 2-3 (shares the string with Overview!G)
$dede.Range("H2").Value = "2016-08-16 18:13:33"
$dede.Range("H3").Value = "2016-08-16 18:13:33"

# Correspond Handback DateTime (K) on de-de, rows 2-3
$dede.Range("K2").Value = "2016-08-16 18:13:53"
$dede.Range("K3").Value = "2016-08-16 18:13:53"
